$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:E to text so numeric-looking strings (e.g. "0.505", "1.00")
# are stored as text, matching the source inlineStr cells instead of being
# auto-converted to numbers by Excel's input parser.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.683.62'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.585.71'
$ws.Range("E3").Value = '  -3.05%  '
$ws.Range("D5").Value = '206.71'
$ws.Range("E5").Value = '  -2.34%  '
$ws.Range("D6").Value = '0.505'
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '22.29'
$ws.Range("E8").Value = '  -4.51%  '
$ws.Range("D9").Value = '0.254'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("E10").Value = '  -3.09%  '
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("D12").Value = '1.811.87'
$ws.Range("E12").Value = '  -3.01%  '
$ws.Range("D13").Value = '1.635.42'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").Value = '3.86'
$ws.Range("E14").Value = '  -3.89%  '
$ws.Range("D15").Value = '0.533'
$ws.Range("E15").Value = '  -5.28%  '
$ws.Range("D16").Value = '27.666.79'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").Value = '63.26'
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("D18").Value = '220.70'
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("D19").Value = '0.0₃0693'
$ws.Range("E19").Value = '  -3.65%  '
$ws.Range("D20").Value = '7.33'
$ws.Range("E20").Value = '  -4.98%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  -5.06%  '
$ws.Range("E23").Value = '  -5.66%  '
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  -5.28%  '
$ws.Range("D25").Value = '153.93'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '6.75'
$ws.Range("E26").Value = '  -2.66%  '
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").Value = '15.12'
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("E29").Value = '  -3.97%  '
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("D33").Value = '1.385.59'
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("E34").Value = '  -4.88%  '
$ws.Range("E35").Value = '  -5.49%  '
$ws.Range("D36").Value = '0.966'
$ws.Range("E36").Value = '  -4.95%  '
$ws.Range("D37").Value = '2.32'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("E38").Value = '  -2.67%  '
$ws.Range("D39").Value = '0.542'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("D40").Value = '0.823'
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").Value = '0.977'
$ws.Range("E42").Value = '  -3.00%  '
$ws.Range("D43").Value = '1.80'
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '63.81'
$ws.Range("E44").Value = '  -3.53%  '
$ws.Range("D45").Value = '2.18'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").Value = '5.24'
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("D47").Value = '1.722.51'
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("D48").Value = '87.92'
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("D49").Value = '0.0₆0101'
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("E50").Value = '  -5.02%  '
$ws.Range("E51").Value = '  -0.94%  '

# Restore the original (default) style on the touched range so only the
# cell values change -- not their formatting/style index.
$ws.Range("D2:E51").Style = "Normal"

